$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: force the date-like text in column A to be stored as plain text
# (matching the other text-date cells already in the sheet) rather than
# letting Excel auto-convert it into a date serial number.
$ws.Cells.Item(34, 1).NumberFormat = "@"
$ws.Cells.Item(34, 1).Value = "07/23/2025"
$ws.Cells.Item(34, 1).ClearFormats()

$ws.Cells.Item(34, 2).Value = 0.0004208400000000022
$ws.Cells.Item(34, 3).Value = 118809.9990495194
$ws.Cells.Item(34, 4).Value = 50
